$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = [double]"1.046867666666667"
$ws.Cells.Item(2, 8).Value = [double]"3.140603"
$ws.Cells.Item(2, 9).Value = [double]"0.000687505225377314"
$ws.Cells.Item(2, 10).Value = [double]"0.000687505225377314"
$ws.Cells.Item(2, 13).Value = [double]"7.487621999999999"
$ws.Cells.Item(2, 14).Value = [double]"22.462866"
$ws.Cells.Item(2, 15).Value = [double]"0.1384395179233961"
$ws.Cells.Item(2, 16).Value = [double]"0.1384395179233961"
$ws.Cells.Item(2, 17).Value = [double]"7.838549372021999"
$ws.Cells.Item(2, 18).Value = [double]"70.546944348198"
$ws.Cells.Item(2, 19).Value = [double]"9.517789197105111E-05"
$ws.Cells.Item(2, 20).Value = [double]"9.517789197105111E-05"
$ws.Cells.Item(3, 7).Value = [double]"1.046867666666667"
$ws.Cells.Item(3, 8).Value = [double]"3.140603"
$ws.Cells.Item(3, 9).Value = [double]"0.000687505225377314"
$ws.Cells.Item(3, 10).Value = [double]"0.000687505225377314"
$ws.Cells.Item(3, 15).Value = [double]"0.5916411627275552"
$ws.Cells.Item(3, 16).Value = [double]"0.5916411627275552"
$ws.Cells.Item(3, 17).Value = [double]"33.499166525028"
$ws.Cells.Item(3, 18).Value = [double]"301.492498725252"
$ws.Cells.Item(3, 19).Value = [double]"0.0004067563909235039"
$ws.Cells.Item(3, 20).Value = [double]"0.0004067563909235039"
$ws.Cells.Item(4, 7).Value = [double]"1.046867666666667"
$ws.Cells.Item(4, 8).Value = [double]"3.140603"
$ws.Cells.Item(4, 9).Value = [double]"0.000687505225377314"
$ws.Cells.Item(4, 10).Value = [double]"0.000687505225377314"
$ws.Cells.Item(4, 13).Value = [double]"14.59882166666667"
$ws.Cells.Item(4, 14).Value = [double]"43.796465"
$ws.Cells.Item(4, 15).Value = [double]"0.2699193193490487"
$ws.Cells.Item(4, 16).Value = [double]"0.2699193193490487"
$ws.Cells.Item(4, 17).Value = [double]"15.28303437426611"
$ws.Cells.Item(4, 18).Value = [double]"137.547309368395"
$ws.Cells.Item(4, 19).Value = [double]"0.0001855709424827589"
$ws.Cells.Item(4, 20).Value = [double]"0.0001855709424827589"
$ws.Cells.Item(5, 8).Value = [double]"4442.55542"
$ws.Cells.Item(5, 9).Value = [double]"0.9725138978974124"
$ws.Cells.Item(5, 10).Value = [double]"0.9725138978974125"
$ws.Cells.Item(5, 13).Value = [double]"7.487621999999999"
$ws.Cells.Item(5, 14).Value = [double]"22.462866"
$ws.Cells.Item(5, 15).Value = [double]"0.1384395179233961"
$ws.Cells.Item(5, 16).Value = [double]"0.1384395179233961"
$ws.Cells.Item(5, 17).Value = [double]"11088.05856633708"
$ws.Cells.Item(5, 18).Value = [double]"99792.5270970337"
$ws.Cells.Item(5, 19).Value = [double]"0.1346343551987206"
$ws.Cells.Item(5, 20).Value = [double]"0.1346343551987206"
$ws.Cells.Item(6, 8).Value = [double]"4442.55542"
$ws.Cells.Item(6, 9).Value = [double]"0.9725138978974124"
$ws.Cells.Item(6, 10).Value = [double]"0.9725138978974125"
$ws.Cells.Item(6, 15).Value = [double]"0.5916411627275552"
$ws.Cells.Item(6, 16).Value = [double]"0.5916411627275552"
$ws.Cells.Item(6, 17).Value = [double]"47386.41076609991"
$ws.Cells.Item(6, 18).Value = [double]"426477.6968948992"
$ws.Cells.Item(6, 19).Value = [double]"0.5753792533207319"
$ws.Cells.Item(6, 20).Value = [double]"0.575379253320732"
$ws.Cells.Item(7, 8).Value = [double]"4442.55542"
$ws.Cells.Item(7, 9).Value = [double]"0.9725138978974124"
$ws.Cells.Item(7, 10).Value = [double]"0.9725138978974125"
$ws.Cells.Item(7, 13).Value = [double]"14.59882166666667"
$ws.Cells.Item(7, 14).Value = [double]"43.796465"
$ws.Cells.Item(7, 15).Value = [double]"0.2699193193490487"
$ws.Cells.Item(7, 16).Value = [double]"0.2699193193490487"
$ws.Cells.Item(7, 17).Value = [double]"21618.69144028781"
$ws.Cells.Item(7, 18).Value = [double]"194568.2229625903"
$ws.Cells.Item(7, 19).Value = [double]"0.2625002893779598"
$ws.Cells.Item(7, 20).Value = [double]"0.2625002893779598"
$ws.Cells.Item(8, 7).Value = [double]"40.80635833333333"
$ws.Cells.Item(8, 8).Value = [double]"122.419075"
$ws.Cells.Item(8, 9).Value = [double]"0.02679859687721029"
$ws.Cells.Item(8, 10).Value = [double]"0.0267985968772103"
$ws.Cells.Item(8, 13).Value = [double]"7.487621999999999"
$ws.Cells.Item(8, 14).Value = [double]"22.462866"
$ws.Cells.Item(8, 15).Value = [double]"0.1384395179233961"
$ws.Cells.Item(8, 16).Value = [double]"0.1384395179233961"
$ws.Cells.Item(8, 17).Value = [double]"305.54258639655"
$ws.Cells.Item(8, 18).Value = [double]"2749.88327756895"
$ws.Cells.Item(8, 19).Value = [double]"0.003709984832704421"
$ws.Cells.Item(8, 20).Value = [double]"0.003709984832704421"
$ws.Cells.Item(9, 7).Value = [double]"40.80635833333333"
$ws.Cells.Item(9, 8).Value = [double]"122.419075"
$ws.Cells.Item(9, 9).Value = [double]"0.02679859687721029"
$ws.Cells.Item(9, 10).Value = [double]"0.0267985968772103"
$ws.Cells.Item(9, 15).Value = [double]"0.5916411627275552"
$ws.Cells.Item(9, 16).Value = [double]"0.5916411627275552"
$ws.Cells.Item(9, 17).Value = [double]"1305.7801254297"
$ws.Cells.Item(9, 18).Value = [double]"11752.0211288673"
$ws.Cells.Item(9, 19).Value = [double]"0.01585515301589973"
$ws.Cells.Item(9, 20).Value = [double]"0.01585515301589973"
$ws.Cells.Item(10, 7).Value = [double]"40.80635833333333"
$ws.Cells.Item(10, 8).Value = [double]"122.419075"
$ws.Cells.Item(10, 9).Value = [double]"0.02679859687721029"
$ws.Cells.Item(10, 10).Value = [double]"0.0267985968772103"
$ws.Cells.Item(10, 13).Value = [double]"14.59882166666667"
$ws.Cells.Item(10, 14).Value = [double]"43.796465"
$ws.Cells.Item(10, 15).Value = [double]"0.2699193193490487"
$ws.Cells.Item(10, 16).Value = [double]"0.2699193193490487"
$ws.Cells.Item(10, 17).Value = [double]"595.7247481744305"
$ws.Cells.Item(10, 18).Value = [double]"5361.522733569875"
$ws.Cells.Item(10, 19).Value = [double]"0.007233459028606145"
$ws.Cells.Item(10, 20).Value = [double]"0.007233459028606146"
